$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RTL freeze criteria / FSM coverage waiver text with the final
# wording used for the RTL Freeze checklist report (adds the note about
# Unreachable States / Unreachable Transitions not being caught).
$ws.Range("F7").Value = "RTL freeze criteria for Dolphin Design but not for OpenHW Group.`nFSM Transition and State are encompassed by Branch, Condition and Statements code coverage.`nBut this doesn't allow to catch Unreachable States and Unreacheable Transitions."

# Leave the selection where the author left it when saving the final version.
$ws.Range("F12").Select() | Out-Null
